$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New columns H and I with headers
$ws.Range("H1").Value = "Image"
$ws.Range("I1").Value = "DOB"

# Row 2 existing columns C/D/E shift to new short test strings, F keeps PhD
$ws.Range("C2").Value = "ash"
$ws.Range("D2").Value = "rai"
$ws.Range("E2").Value = "bac"
$ws.Range("F2").Value = "PhD"

# New row2 values for H and I
$ws.Range("H2").Value = "iii"
$ws.Range("I2").Value = "28-04-2028"
$ws.Range("I2").NumberFormat = "@"

$ws.Range("I2").Select()
